function ColLetter($n) {
  $letter = ""
  while ($n -gt 0) {
    $rem = ($n - 1) % 26
    $letter = [char](65 + $rem) + $letter
    $n = [int](($n - $rem - 1) / 26)
  }
  return $letter
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table grew from 9 columns (A:I) x 9 rows to 11 columns (A:K) x 13 rows,
# with two new columns ("ownTeam", "oppTeam") inserted and the rows reordered/
# extended with additional match data. Clear everything first, then rewrite
# the whole block fresh so stale cells/formatting from the old 9x9 layout do
# not linger outside the new A1:K13 range.
$ws.Cells.Clear()

$lastCol = ColLetter 11
$addr = "A1:" + $lastCol + "13"
$rng = $ws.Range($addr)

# Force text storage so numeric-looking values keep their exact original
# formatting (e.g. "120.00" stays "120.00" instead of becoming the number 120,
# and "0" stays text "0") -- matching the source data which stores every
# column, including numeric ones, as text.
$rng.NumberFormat = "@"

$data = New-Object 'object[,]' 13,11
$data[0,0] = 'venue'
$data[0,1] = 'date'
$data[0,2] = 'result'
$data[0,3] = 'ownTeam'
$data[0,4] = 'oppTeam'
$data[0,5] = 'batsman'
$data[0,6] = 'totalRuns'
$data[0,7] = 'totalBalls'
$data[0,8] = 'total4s'
$data[0,9] = 'total6s'
$data[0,10] = 'sr'
$data[1,0] = ' Sharjah'
$data[1,1] = ' November 03 2020'
$data[1,2] = 'Sunrisers won by 10 wickets (with 17 balls remaining)'
$data[1,3] = 'Mumbai Indians'
$data[1,4] = 'Sunrisers Hyderabad'
$data[1,5] = 'Rohit Sharma (c)'
$data[1,6] = '4'
$data[1,7] = '7'
$data[1,8] = '0'
$data[1,9] = '0'
$data[1,10] = '57.14'
$data[2,0] = ' Dubai (DSC)'
$data[2,1] = ' November 05 2020'
$data[2,2] = 'Mumbai won by 57 runs'
$data[2,3] = 'Mumbai Indians'
$data[2,4] = 'Delhi Capitals'
$data[2,5] = 'Rohit Sharma (c)'
$data[2,6] = '0'
$data[2,7] = '1'
$data[2,8] = '0'
$data[2,9] = '0'
$data[2,10] = '0.00'
$data[3,0] = ' Dubai (DSC)'
$data[3,1] = ' November 10 2020'
$data[3,2] = 'Mumbai won by 5 wickets (with 8 balls remaining)'
$data[3,3] = 'Mumbai Indians'
$data[3,4] = 'Delhi Capitals'
$data[3,5] = 'Rohit Sharma (c)'
$data[3,6] = '68'
$data[3,7] = '51'
$data[3,8] = '5'
$data[3,9] = '4'
$data[3,10] = '133.33'
$data[4,0] = ' Abu Dhabi'
$data[4,1] = ' October 16 2020'
$data[4,2] = 'Mumbai won by 8 wickets (with 19 balls remaining)'
$data[4,3] = 'Mumbai Indians'
$data[4,4] = 'Kolkata Knight Riders'
$data[4,5] = 'Rohit Sharma (c)'
$data[4,6] = '35'
$data[4,7] = '36'
$data[4,8] = '5'
$data[4,9] = '1'
$data[4,10] = '97.22'
$data[5,0] = ' Abu Dhabi'
$data[5,1] = ' October 01 2020'
$data[5,2] = 'Mumbai won by 48 runs'
$data[5,3] = 'Mumbai Indians'
$data[5,4] = 'Kings XI Punjab'
$data[5,5] = 'Rohit Sharma (c)'
$data[5,6] = '70'
$data[5,7] = '45'
$data[5,8] = '8'
$data[5,9] = '3'
$data[5,10] = '155.55'
$data[6,0] = ' Abu Dhabi'
$data[6,1] = ' September 19 2020'
$data[6,2] = 'Super Kings won by 5 wickets (with 4 balls remaining)'
$data[6,3] = 'Mumbai Indians'
$data[6,4] = 'Chennai Super Kings'
$data[6,5] = 'Rohit Sharma (c)'
$data[6,6] = '12'
$data[6,7] = '10'
$data[6,8] = '2'
$data[6,9] = '0'
$data[6,10] = '120.00'
$data[7,0] = ' Abu Dhabi'
$data[7,1] = ' October 06 2020'
$data[7,2] = 'Mumbai won by 57 runs'
$data[7,3] = 'Mumbai Indians'
$data[7,4] = 'Rajasthan Royals'
$data[7,5] = 'Rohit Sharma (c)'
$data[7,6] = '35'
$data[7,7] = '23'
$data[7,8] = '2'
$data[7,9] = '3'
$data[7,10] = '152.17'
$data[8,0] = ' Dubai (DSC)'
$data[8,1] = ' October 18 2020'
$data[8,2] = 'Match tied (Kings XI won the one-over eliminator)'
$data[8,3] = 'Mumbai Indians'
$data[8,4] = 'Kings XI Punjab'
$data[8,5] = 'Rohit Sharma (c)'
$data[8,6] = '9'
$data[8,7] = '8'
$data[8,8] = '2'
$data[8,9] = '0'
$data[8,10] = '112.50'
$data[9,0] = ' Dubai (DSC)'
$data[9,1] = ' September 28 2020'
$data[9,2] = 'Match tied (RCB won the one-over eliminator)'
$data[9,3] = 'Mumbai Indians'
$data[9,4] = 'Royal Challengers Bangalore'
$data[9,5] = 'Rohit Sharma (c)'
$data[9,6] = '8'
$data[9,7] = '8'
$data[9,8] = '0'
$data[9,9] = '1'
$data[9,10] = '100.00'
$data[10,0] = ' Abu Dhabi'
$data[10,1] = ' September 23 2020'
$data[10,2] = 'Mumbai won by 49 runs'
$data[10,3] = 'Mumbai Indians'
$data[10,4] = 'Kolkata Knight Riders'
$data[10,5] = 'Rohit Sharma (c)'
$data[10,6] = '80'
$data[10,7] = '54'
$data[10,8] = '3'
$data[10,9] = '6'
$data[10,10] = '148.14'
$data[11,0] = ' Abu Dhabi'
$data[11,1] = ' October 11 2020'
$data[11,2] = 'Mumbai won by 5 wickets (with 2 balls remaining)'
$data[11,3] = 'Mumbai Indians'
$data[11,4] = 'Delhi Capitals'
$data[11,5] = 'Rohit Sharma (c)'
$data[11,6] = '5'
$data[11,7] = '12'
$data[11,8] = '0'
$data[11,9] = '0'
$data[11,10] = '41.66'
$data[12,0] = ' Sharjah'
$data[12,1] = ' October 04 2020'
$data[12,2] = 'Mumbai won by 34 runs'
$data[12,3] = 'Mumbai Indians'
$data[12,4] = 'Sunrisers Hyderabad'
$data[12,5] = 'Rohit Sharma (c)'
$data[12,6] = '6'
$data[12,7] = '5'
$data[12,8] = '0'
$data[12,9] = '1'
$data[12,10] = '120.00'

$rng.Value = $data

Write-Host "Applied data update to sheet1 ($addr)"
